# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets
# to reflect newly generated output (commit: 456a3b4).

$wb = $excel.ActiveWorkbook

# Map of worksheet name -> { row -> newValue }
$updates = @{
    "展览" = @{
        2  = 12552
        4  = 2036
        5  = 267
        6  = 387
        8  = 12520
        9  = 3084
        10 = 533
        11 = 3
        16 = 648
        17 = 2837
        18 = 6082
    }
    "全部类型" = @{
        2  = 12552
        4  = 2036
        5  = 267
        7  = 387
        9  = 12520
        10 = 3084
        11 = 533
        12 = 3
        17 = 648
        18 = 2837
        20 = 6082
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $newValue = $rows[$row]
        $ws.Range("F$row").Value = $newValue
    }
}
